$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "523788",
    "630822",
    "414931",
    "647554",
    "704511",
    "438543",
    "319834",
    "376791",
    "982415",
    "139373",
    "676223",
    "607590",
    "244442",
    "928550",
    "584098",
    "239646",
    "873678",
    "353560",
    "410517",
    "370292",
    "301659",
    "211358",
    "318392",
    "102501",
    "385201",
    "843566",
    "452009",
    "236119",
    "920227",
    "654260"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 4)
    $cell.Select()
    $cell.Value = $values[$i]
}
